$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data: two string columns plus a date and a date+time value.
$ws.Range("A10").Value = "data8"
$ws.Range("B10").Value = "row with two dates"

$ws.Range("C10").Value = 40909
$ws.Range("C10").NumberFormat = "MM/DD/YY"

$ws.Range("D10").Value = 40953.0930555556
$ws.Range("D10").NumberFormat = "MM/DD/YYYY\ HH:MM:SS"

# Leave the selection where it ends up after the edit, matching the author's
# recorded cursor position.
[void]$ws.Range("A10").Select()
